$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The original F3 formula was a broken/mistyped array formula:
#   F3: {=F2+F18E3}  -> #NAME? error, which then cascaded #NAME? errors
#   down through F4:F12 and G4:G12 (their dependents).
# Correct it back to a normal (non-array) formula: F2 + E3.
$ws.Range("F3").Formula = "=F2+E3"

# F4 continues the running total pattern (not part of the F5:F12 shared group).
$ws.Range("F4").Formula = "=F3+E4"

# F5:F12 form a shared-formula block continuing the same running-total pattern.
$ws.Range("F5:F12").Formula = "=F4+E5"

# G4:G12 recompute automatically now that the F column errors are gone
# (their formulas were already correct, just poisoned by the F-column errors).

# Restore the sheet view: zoomed to 125% with B27 as the active selection.
$excel.ActiveWindow.Zoom = 125
$ws.Range("B27").Select()

$wb.Save()
